$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper style-source cells that keep their style (text=14, whole-number=15,
# one-decimal-percent=16) throughout this edit, used to "re-stamp" the
# correct cellXf after a value/type change flips a cell between text and
# number (PasteSpecial formats-only keeps the numeric value intact).
$TEXT_STYLE_SRC = "A16"   # s="14" (General / text)
$NUM_STYLE_SRC  = "I16"   # s="15" (#,##0)
$PCT_STYLE_SRC  = "K16"   # s="16" (#,##0.0)

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $ws.Range($TEXT_STYLE_SRC).Copy()
    $r.PasteSpecial(-4122)
}

function Set-NumCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.Value = $val
    $ws.Range($NUM_STYLE_SRC).Copy()
    $r.PasteSpecial(-4122)
}

function Set-PctCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.Value = $val
    $ws.Range($PCT_STYLE_SRC).Copy()
    $r.PasteSpecial(-4122)
}

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/25/2023  Through  12/31/2023"

# --- Row 14 ---
Set-TextCell "F14" "0"
$ws.Range("N14").Value = -90.909090909090

# --- Row 15 ---
Set-TextCell "D15" "0"
Set-TextCell "E15" "***.*"

# --- Row 16 ---
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 800
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 107.142857142857
$ws.Range("I16").Value = 390
$ws.Range("J16").Value = 316
$ws.Range("K16").Value = 23.417721518987
$ws.Range("L16").Value = 54.150197628458
$ws.Range("M16").Value = -15.217391304347
$ws.Range("N16").Value = -74.476439790575

# --- Row 17 ---
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = 37.5
$ws.Range("I17").Value = 472
$ws.Range("J17").Value = 376
$ws.Range("K17").Value = 25.531914893617
$ws.Range("L17").Value = 47.962382445141
$ws.Range("M17").Value = 25.531914893617
$ws.Range("N17").Value = -44.142011834319

# --- Row 18 ---
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 248
$ws.Range("J18").Value = 299
$ws.Range("K18").Value = -17.056856187291
$ws.Range("L18").Value = -10.144927536231
$ws.Range("M18").Value = -48.225469728601
$ws.Range("N18").Value = -82.106782106782

# --- Row 19 ---
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -30.769230769230
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = 20
$ws.Range("I19").Value = 716
$ws.Range("J19").Value = 702
$ws.Range("K19").Value = 1.994301994301
$ws.Range("L19").Value = 21.768707482993
$ws.Range("M19").Value = 138.666666666667
$ws.Range("N19").Value = 26.501766784452

# --- Row 20 ---
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 500
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 125
$ws.Range("I20").Value = 234
$ws.Range("J20").Value = 239
$ws.Range("K20").Value = -2.092050209205
$ws.Range("L20").Value = 35.260115606936
$ws.Range("M20").Value = 32.954545454545
$ws.Range("N20").Value = -75.159235668789

# --- Row 21 ---
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 29.166666666666
$ws.Range("F21").Value = 162
$ws.Range("G21").Value = 116
$ws.Range("H21").Value = 39.655172413793
$ws.Range("I21").Value = 2086
$ws.Range("J21").Value = 1975
$ws.Range("K21").Value = 5.620253164556
$ws.Range("L21").Value = 26.885644768856
$ws.Range("M21").Value = 13.492927094668
$ws.Range("N21").Value = -61.176251628512

# --- Row 22 ---
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 28
$ws.Range("J22").Value = 19
$ws.Range("K22").Value = 47.368421052631
$ws.Range("L22").Value = 21.739130434782
$ws.Range("M22").Value = 7.692307692307

# --- Row 23 ---
Set-TextCell "C23" "0"
Set-NumCell "D23" 1
Set-PctCell "E23" -100
$ws.Range("J23").Value = 29
$ws.Range("K23").Value = 31.034482758620

# --- Row 24 ---
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 28.571428571428
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 75
$ws.Range("H24").Value = 20
$ws.Range("I24").Value = 1007
$ws.Range("J24").Value = 1056
$ws.Range("K24").Value = -4.640151515151
$ws.Range("L24").Value = 3.282051282051
$ws.Range("M24").Value = 33.025099075297

# --- Row 25 ---
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 142.857142857143
$ws.Range("F25").Value = 53
$ws.Range("G25").Value = 47
$ws.Range("H25").Value = 12.765957446808
$ws.Range("I25").Value = 692
$ws.Range("J25").Value = 587
$ws.Range("K25").Value = 17.887563884156
$ws.Range("L25").Value = 31.309297912713
$ws.Range("M25").Value = -10.594315245478

# --- Row 26 ---
Set-TextCell "C26" "0"
Set-TextCell "D26" "0"
Set-TextCell "E26" "***.*"

# --- Row 27 ---
Set-TextCell "C27" "0"
Set-NumCell "D27" 1
Set-PctCell "E27" -100
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 350
$ws.Range("J27").Value = 63
$ws.Range("K27").Value = 7.936507936507

# --- Row 28 ---
$ws.Range("G28").Value = 2
$ws.Range("N28").Value = -93.697478991596

# --- Row 29 ---
$ws.Range("G29").Value = 1
$ws.Range("N29").Value = -93.577981651376

# --- Row 30 ---
Set-NumCell "F30" 1
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 8
$ws.Range("K30").Value = 33.333333333333
$ws.Range("L30").Value = 166.666666666667

$wb.Save()
